$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (translated to Spanish)
$ws.Range("A1").Value = "Número de cédula"
$ws.Range("D1").Value = "Número de celular"

# Remove the hyperlinks tied to the email addresses in column B before
# deleting the rows that hosted them.
$ws.Hyperlinks.Delete()

# Drop the sample data rows (2-5), leaving only the header row.
$ws.Range("A2:F5").EntireRow.Delete()

# Put the active selection back on the first data-entry cell of the header row.
$ws.Range("D1").Select()
